# relevant_variables.xlsx — mark a handful of "Categorical" rows as "Target"
# and un-filter the sheet (show all rows, clear the Type=Categorical filter).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Remove the autofilter criterion on column C ("Type") so every row is
#     shown again (un-hides the rows that the filter had hidden, while
#     keeping the autofilter dropdown buttons themselves in place). ---
$ws.Range("A1:D44").AutoFilter(3) | Out-Null

# --- Relabel the target variable rows' Type from "Categorical" to "Target" ---
$targetRows = 3, 4, 5, 13, 14, 15, 41
foreach ($r in $targetRows) {
    $ws.Cells.Item($r, 3).Value2 = "Target"
}

# Row 41 (SOCSCLPAR_A) becomes a highlighted row like the other Target rows
# (rows 3-5, 13-15 already carry the yellow "customFormat" fill).
$ws.Range("A41:C41").Interior.Color = $ws.Range("A3").Interior.Color

# --- Update the active selection to reflect where the editor left off ---
$ws.Range("C42").Select() | Out-Null
